# Auto-generated script applying scheduled-runner market data updates
# to the Sophia_Profits workbook (per-sheet Leve profit columns H,I,J,K,L,M,N).
$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2778.2144
$ws.Range("I15").Value = 2778.2144
$ws.Range("K15").Value = 8334.643199999999
$ws.Range("M15").Value = -8165.643199999999

# ALC row 32
$ws.Range("H32").Value = 1289
$ws.Range("I32").Value = 1482.3334
$ws.Range("J32").Value = 999
$ws.Range("K32").Value = 1482.3334
$ws.Range("L32").Value = 999
$ws.Range("M32").Value = -1156.3334
$ws.Range("N32").Value = -1651

# ALC row 55
$ws.Range("H55").Value = 255.55556
$ws.Range("I55").Value = 255.55556
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 255.55556
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -41.55556000000001
$ws.Range("N55").ClearContents()

# ALC row 98
$ws.Range("H98").Value = 7666.6665
$ws.Range("I98").Value = 2000
$ws.Range("K98").Value = 2000
$ws.Range("M98").Value = -502

# ALC row 115
$ws.Range("H115").Value = 3107.7144
$ws.Range("I115").Value = 1125.6666
$ws.Range("J115").Value = 15000
$ws.Range("K115").Value = 3376.9998
$ws.Range("L115").Value = 45000
$ws.Range("M115").Value = -1809.9998
$ws.Range("N115").Value = -48134

# ALC row 116
$ws.Range("H116").Value = 7199.5
$ws.Range("I116").Value = 9900
$ws.Range("K116").Value = 9900
$ws.Range("M116").Value = -6458

# ALC row 122
$ws.Range("H122").Value = 7666.6665
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

# ALC row 131
$ws.Range("H131").Value = 3765.8333
$ws.Range("J131").Value = 3560
$ws.Range("L131").Value = 10680
$ws.Range("N131").Value = -20760

# ALC row 132
$ws.Range("H132").Value = 1151.1875
$ws.Range("I132").Value = 1185.3334
$ws.Range("K132").Value = 3556.0002
$ws.Range("M132").Value = -1026.0002

# ALC row 137
$ws.Range("H137").Value = 980.6
$ws.Range("I137").Value = 991.8889
$ws.Range("K137").Value = 2975.6667
$ws.Range("M137").Value = -425.6667000000002

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2603.3333
$ws.Range("I45").Value = 3047.2856
$ws.Range("K45").Value = 3047.2856
$ws.Range("M45").Value = -2670.2856

# ARM row 61
$ws.Range("H61").Value = 1817.6666
$ws.Range("I61").Value = 1476.5
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 1476.5
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -1264.5
$ws.Range("N61").Value = -2924

# ARM row 110
$ws.Range("H110").Value = 3790.6924
$ws.Range("I110").Value = 2139.5454
$ws.Range("K110").Value = 2139.5454
$ws.Range("M110").Value = -94.54539999999997

# ARM row 136
$ws.Range("H136").Value = 1817.6666
$ws.Range("I136").Value = 1476.5
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 4429.5
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -1879.5
$ws.Range("N136").Value = -12600

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2164.5
$ws.Range("I99").Value = 2164.5
$ws.Range("K99").Value = 2164.5
$ws.Range("M99").Value = -666.5

# BSM row 105
$ws.Range("H105").Value = 3998.25
$ws.Range("I105").Value = 3998.25
$ws.Range("K105").Value = 3998.25
$ws.Range("M105").Value = -2251.25

# BSM row 134
$ws.Range("H134").Value = 2360.889
$ws.Range("I134").Value = 2349.7144
$ws.Range("J134").Value = 2400
$ws.Range("K134").Value = 7049.1432
$ws.Range("L134").Value = 7200
$ws.Range("M134").Value = -4514.1432
$ws.Range("N134").Value = -12270

# BSM row 141
$ws.Range("H141").Value = 179996.5
$ws.Range("J141").Value = 199997
$ws.Range("L141").Value = 199997
$ws.Range("N141").Value = -210357

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1830.8125
$ws.Range("I31").Value = 1269.6
$ws.Range("J31").Value = 2766.1667
$ws.Range("K31").Value = 1269.6
$ws.Range("L31").Value = 2766.1667
$ws.Range("M31").Value = -974.5999999999999
$ws.Range("N31").Value = -3356.1667

# CRP row 34
$ws.Range("H34").Value = 1830.8125
$ws.Range("I34").Value = 1269.6
$ws.Range("J34").Value = 2766.1667
$ws.Range("K34").Value = 1269.6
$ws.Range("L34").Value = 2766.1667
$ws.Range("M34").Value = -1067.6
$ws.Range("N34").Value = -3170.1667

# CRP row 132
$ws.Range("H132").Value = 2599.6155
$ws.Range("I132").Value = 1846.4445
$ws.Range("K132").Value = 5539.333500000001
$ws.Range("M132").Value = -3009.333500000001

# CRP row 141
$ws.Range("H141").Value = 149994
$ws.Range("J141").Value = 149994
$ws.Range("L141").Value = 149994
$ws.Range("N141").Value = -160354

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 204597.8
$ws.Range("J81").Value = 338330
$ws.Range("L81").Value = 1014990
$ws.Range("N81").Value = -1017236

# CUL row 84
$ws.Range("H84").Value = 204597.8
$ws.Range("J84").Value = 338330
$ws.Range("L84").Value = 3044970
$ws.Range("N84").Value = -3056202

# CUL row 87
$ws.Range("H87").Value = 8134.143
$ws.Range("I87").Value = 6823.1665
$ws.Range("K87").Value = 20469.4995
$ws.Range("M87").Value = -19221.4995

# CUL row 90
$ws.Range("H90").Value = 8134.143
$ws.Range("I90").Value = 6823.1665
$ws.Range("K90").Value = 61408.4985
$ws.Range("M90").Value = -55168.4985

# CUL row 129
$ws.Range("H129").Value = 1968.4445
$ws.Range("I129").Value = 1289
$ws.Range("K129").Value = 3867
$ws.Range("M129").Value = 1133

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3195.875
$ws.Range("I113").Value = 3195.875
$ws.Range("K113").Value = 3195.875
$ws.Range("M113").Value = -1025.875

# GSM row 122
$ws.Range("H122").Value = 3703.9
$ws.Range("I122").Value = 3904.3333
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 11712.9999
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -9262.999899999999
$ws.Range("N122").Value = -10600

# GSM row 132
$ws.Range("H132").Value = 4632.5
$ws.Range("J132").Value = 4599.3335
$ws.Range("L132").Value = 13798.0005
$ws.Range("N132").Value = -18858.0005

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7293.9287
$ws.Range("I22").Value = 5895
$ws.Range("J22").Value = 8692.857
$ws.Range("K22").Value = 5895
$ws.Range("L22").Value = 8692.857
$ws.Range("M22").Value = -5600
$ws.Range("N22").Value = -9282.857

# LTW row 27
$ws.Range("H27").Value = 7293.9287
$ws.Range("I27").Value = 5895
$ws.Range("J27").Value = 8692.857
$ws.Range("K27").Value = 5895
$ws.Range("L27").Value = 8692.857
$ws.Range("M27").Value = -5788
$ws.Range("N27").Value = -8906.857

# LTW row 40
$ws.Range("H40").Value = 3150.8
$ws.Range("I40").Value = 3333
$ws.Range("J40").Value = 2877.5
$ws.Range("K40").Value = 3333
$ws.Range("L40").Value = 2877.5
$ws.Range("M40").Value = -3197
$ws.Range("N40").Value = -3149.5

# LTW row 82
$ws.Range("H82").Value = 3426.4285
$ws.Range("I82").Value = 3600
$ws.Range("K82").Value = 3600
$ws.Range("M82").Value = -3239

# LTW row 85
$ws.Range("H85").Value = 3426.4285
$ws.Range("I85").Value = 3600
$ws.Range("K85").Value = 3600
$ws.Range("M85").Value = -2352

# LTW row 136
$ws.Range("H136").Value = 2907
$ws.Range("I136").Value = 2907
$ws.Range("K136").Value = 8721
$ws.Range("M136").Value = -6171

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5478.857
$ws.Range("I122").Value = 5478.857
$ws.Range("K122").Value = 16436.571
$ws.Range("M122").Value = -13986.571

# WVR row 132
$ws.Range("H132").Value = 2737
$ws.Range("I132").Value = 2474.5
$ws.Range("K132").Value = 7423.5
$ws.Range("M132").Value = -4893.5

# WVR row 136
$ws.Range("H136").Value = 3625.5334
$ws.Range("I136").Value = 3655.1
$ws.Range("J136").Value = 3566.4
$ws.Range("K136").Value = 10965.3
$ws.Range("L136").Value = 10699.2
$ws.Range("M136").Value = -8415.299999999999
$ws.Range("N136").Value = -15799.2
